# --- edit.ps1 -------------------------------------------------------------
# New crime data collected: update the volume/date header strings and the
# Week to Date / 28 Day / Year to Date / 2 Year crime-complaint figures for
# rows 15-30 (Murder..Hate Crimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text: volume number and the two report-covering dates.
#    Both shared strings are rich-text (multiple runs); only the leaf run
#    holding the actual number/date needs to change, so we patch the exact
#    character ranges via .Characters(start, length).Text rather than
#    replacing the whole cell value.
# ---------------------------------------------------------------------------
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "46"   # "...Number  45" -> "...Number  46"

$weekCell = $ws.Range("C9")
# Replace the second (later) date first so the first replacement's offsets
# are not shifted by the differing string length ("11/13/2022" -> "11/20/2022").
$weekCell.Characters(47, 10).Text = "11/20/2022"
$weekCell.Characters(27, 9).Text = "11/14/2022"

# ---------------------------------------------------------------------------
# 2. Helper functions for the crime-complaint grid (rows 15-30, cols C:N).
# ---------------------------------------------------------------------------
function Set-NumCell($addr, $num, $fmt) {
    $c = $ws.Range($addr)
    $c.Value = $num
    $c.NumberFormat = $fmt
}

function Set-TextCell($addr, $text, $donorAddr) {
    # Force the cell to hold a text value (so "0" / "***.*" are not
    # reinterpreted as numbers), then restore the donor cell's original
    # (non quote-prefixed) style so the saved XF matches the untouched cells.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $c.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 3. Cells that flip between the text placeholder ("0" / "***.*") and a real
#    number (or vice versa) - these need both value AND style/type fixed up.
# ---------------------------------------------------------------------------
Set-NumCell "D15" 1 "#,##0"
Set-NumCell "E15" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "D18" 2 "#,##0"
Set-NumCell "E18" 50 "#,##0.0;""-""#,##0.0"
Set-TextCell "C22" "0" "C14"
Set-TextCell "C23" "0" "C14"
Set-NumCell "D26" 1 "#,##0"
Set-NumCell "E26" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "C27" 1 "#,##0"
Set-TextCell "D27" "0" "C14"
Set-TextCell "E27" "***.*" "E14"
Set-TextCell "C30" "0" "C14"

# ---------------------------------------------------------------------------
# 4. Remaining cells: same style/type as before, only the numeric value
#    itself changes with the new weekly crime totals.
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -7.692307692307
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 192
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = 24.675324675324
$ws.Range("L16").Value = 82.857142857142
$ws.Range("M16").Value = -10.280373831775
$ws.Range("N16").Value = -53.623188405797
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 44.444444444444
$ws.Range("I17").Value = 237
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = 16.176470588235
$ws.Range("L17").Value = 11.267605633802
$ws.Range("M17").Value = 43.636363636363
$ws.Range("N17").Value = 15.04854368932
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 93
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = -6.060606060606
$ws.Range("L18").Value = 3.333333333333
$ws.Range("M18").Value = -66.304347826087
$ws.Range("N18").Value = -87.398373983739
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -17.460317460317
$ws.Range("I19").Value = 514
$ws.Range("J19").Value = 430
$ws.Range("K19").Value = 19.53488372093
$ws.Range("L19").Value = 32.8165374677
$ws.Range("M19").Value = 18.706697459584
$ws.Range("N19").Value = 31.457800511509
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 269
$ws.Range("J20").Value = 278
$ws.Range("K20").Value = -3.237410071942
$ws.Range("L20").Value = 94.927536231884
$ws.Range("M20").Value = 43.085106382978
$ws.Range("N20").Value = -85.94566353187
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = -2.941176470588
$ws.Range("I21").Value = 1321
$ws.Range("J21").Value = 1183
$ws.Range("K21").Value = 11.665257819104
$ws.Range("L21").Value = 39.052631578947
$ws.Range("M21").Value = 2.086553323029
$ws.Range("N21").Value = -64.258658008658
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -20
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = -30.232558139534
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = -38.775510204081
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -17.857142857142
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -8.823529411764
$ws.Range("I24").Value = 1027
$ws.Range("J24").Value = 883
$ws.Range("K24").Value = 16.308040770101
$ws.Range("L24").Value = 26.01226993865
$ws.Range("M24").Value = -21.90114068441
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 40.74074074074
$ws.Range("I25").Value = 407
$ws.Range("J25").Value = 387
$ws.Range("K25").Value = 5.16795865633
$ws.Range("L25").Value = 17.29106628242
$ws.Range("M25").Value = 2.777777777777
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = 52.941176470588
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 46
$ws.Range("K27").Value = 21.052631578947
$ws.Range("L27").Value = 130
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -7.142857142857
$ws.Range("L28").Value = 62.5
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -18.181818181818
$ws.Range("L29").Value = 12.5
